# Update "想去人数" (interest counts) in column F for the duplicated event rows
# that appear on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 10, 12, 13, 16
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F10").Value = 6570
$wsExhibit.Range("F12").Value = 359
$wsExhibit.Range("F13").Value = 2783
$wsExhibit.Range("F16").Value = 254

# Sheet "全部类型": rows 13, 16, 17, 20 (same events, duplicated listing)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F13").Value = 6570
$wsAll.Range("F16").Value = 359
$wsAll.Range("F17").Value = 2783
$wsAll.Range("F20").Value = 254
